$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.142.09'

# Row 3
$ws.Range("D3").Value = '1.680.36'
$ws.Range("E3").Value = '  +0.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.31'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.24%  '

# Row 6
$ws.Range("E6").Value = '  +0.46%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("E8").Value = '  +2.22%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.42'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +5.33%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0623'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.63%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0888'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.16%  '

# Row 12
$ws.Range("D12").Value = '1.917.38'
$ws.Range("E12").Value = '  +0.38%  '

# Row 13
$ws.Range("D13").Value = '1.666.05'
$ws.Range("E13").Value = '  -0.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.60%  '

# Row 15
$ws.Range("E15").Value = '  +2.16%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.28'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.01%  '

# Row 17
$ws.Range("D17").Value = '27.139.18'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '239.19'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.47%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.07'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.33%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0743'
$ws.Range("E20").Value = '  +1.36%  '

# Row 21
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("E22").Value = '  +1.96%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.47'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.99%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.35%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.26'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.71%  '

# Row 27
$ws.Range("E27").Value = '  +2.40%  '

# Row 28
$ws.Range("E28").Value = '  +0.21%  '

# Row 29
$ws.Range("E29").Value = '  -0.11%  '

# Row 30
$ws.Range("E30").Value = '  +0.34%  '

# Row 31
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("D32").Value = '1.561.25'
$ws.Range("E32").Value = '  +5.39%  '

# Row 33
$ws.Range("E33").Value = '  +1.63%  '

# Row 34
$ws.Range("E34").Value = '  +2.88%  '

# Row 35
$ws.Range("E35").Value = '  +0.16%  '

# Row 36
$ws.Range("E36").Value = '  +2.31%  '

# Row 37
$ws.Range("E37").Value = '  -1.19%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.932'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.98%  '

# Row 39
$ws.Range("E39").Value = '  +2.15%  '

# Row 40
$ws.Range("E40").Value = '  +2.43%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.21'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.13%  '

# Row 42
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.58'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.29%  '

# Row 44
$ws.Range("E44").Value = '  -2.23%  '

# Row 45
$ws.Range("D45").Value = '1.825.58'
$ws.Range("E45").Value = '  +0.61%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.787'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.50%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.75'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.19%  '

# Row 48
$ws.Range("E48").Value = '  +3.49%  '

# Row 49
$ws.Range("E49").Value = '  +1.44%  '

# Row 50
$ws.Range("E50").Value = '  +2.22%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.12'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.22%  '
